$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- The model was retrained: fill in the previously-empty Accuracy/Loss
# --- values for the last three rows (MobileNetV2(alpha=1.0), MobileNet
# --- (alpha=1.0), EfficientNet-B0) and apply the same numeric format as
# --- the rest of the table to the whole Accuracy/Loss column range.
$ws.Range("B9").Value = 0.29797556996345498
$ws.Range("C9").Value = 0.90760868787765503
$ws.Range("B10").Value = 0.41326564550399703
$ws.Range("C10").Value = 0.85778987407684304
$ws.Range("B11").Value = 0.31138986349105802
$ws.Range("C11").Value = 0.91123187541961603

$ws.Range("B2:C11").NumberFormat = "0.000_ "

# --- Column widths were narrowed for B/C, and D (row count, always 224) was hidden.
$ws.Columns("B").ColumnWidth = 6.5
$ws.Columns("C").ColumnWidth = 8.2857142857142857
$ws.Columns("D").ColumnWidth = 5.2857142857142857
$ws.Columns("D").Hidden = $true

# --- Selection cursor ended on F11 after the edits.
$ws.Range("F11").Select()
